$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.476000308990479
$ws.Range("B1").Value = 3.638644218444824
$ws.Range("C1").Value = 6.069168090820312
$ws.Range("D1").Value = 1.485948204994202
$ws.Range("E1").Value = 0.8701711297035217
